$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.326.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6980"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08120"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.878.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7283"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.210"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.609.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.912"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007756"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.155.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.641"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.063"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.942"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.402"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.507"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.422"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.060"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05248"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.200"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7206"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.27%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.665"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01870"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.711"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8854"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.911"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.044.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.12%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.276"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.038.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  -5.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.263"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.30%  "
